$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '68.086.05'
$ws.Range("E2").Value = '  +0.70%  '

$ws.Range("D3").Value = '3.337.81'
$ws.Range("E3").Value = '  +1.02%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.07%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '584.15'
$ws.Range("E5").Value = '  +0.74%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '176.03'
$ws.Range("E6").Value = '  +0.97%  '

$ws.Range("E7").Value = '  -0.02%  '

$ws.Range("E8").Value = '  +2.22%  '

$ws.Range("E9").Value = '  +5.01%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.582'
$ws.Range("E10").Value = '  +1.40%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '47.59'
$ws.Range("E11").Value = '  +5.24%  '

$ws.Range("E12").Value = '  +2.11%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '699.67'
$ws.Range("E13").Value = '  +4.77%  '

$ws.Range("D14").Value = '3.885.83'
$ws.Range("E14").Value = '  +1.18%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '8.41'
$ws.Range("E15").Value = '  +0.74%  '

$ws.Range("D16").Value = '68.103.30'
$ws.Range("E16").Value = '  +0.66%  '

$ws.Range("E17").Value = '  +0.94%  '

$ws.Range("D18").Value = '3.339.01'
$ws.Range("E18").Value = '  +1.06%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '17.48'
$ws.Range("E19").Value = '  +0.65%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.15'
$ws.Range("E20").Value = '  +3.15%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.896'
$ws.Range("E21").Value = '  +1.33%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.39'
$ws.Range("E22").Value = '  +0.40%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '17.02'

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '101.15'
$ws.Range("E24").Value = '  +4.14%  '

$ws.Range("E25").Value = '  +2.46%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.69'
$ws.Range("E26").Value = '  +0.92%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.44'
$ws.Range("E27").Value = '  +3.02%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '33.25'
$ws.Range("E28").Value = '  +0.82%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '8.54'
$ws.Range("E29").Value = '  +2.25%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.01'
$ws.Range("E30").Value = '  -0.97%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '571.76'
$ws.Range("E31").Value = '  -3.23%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '11.02'
$ws.Range("E32").Value = '  +1.01%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.105'
$ws.Range("E33").Value = '  +2.34%  '

$ws.Range("D34").Value = '3.761.58'
$ws.Range("E34").Value = '  +1.06%  '

$ws.Range("E35").Value = '  +0.05%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '56.97'
$ws.Range("E36").Value = '  +3.19%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.35'
$ws.Range("E37").Value = '  +0.95%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '35.45'
$ws.Range("E38").Value = '  +10.35%  '

$ws.Range("E39").Value = '  +1.79%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.15'
$ws.Range("E40").Value = '  +1.88%  '

$ws.Range("E41").Value = '  +0.43%  '

$ws.Range("D42").Value = '0.0₃0678'
$ws.Range("E42").Value = '  +2.46%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.335'
$ws.Range("E43").Value = '  +1.79%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '3.31'
$ws.Range("E44").Value = '  +0.70%  '

$ws.Range("E45").Value = '  +0.84%  '

$ws.Range("E46").Value = '  +1.71%  '

$ws.Range("E47").Value = '  +1.72%  '

$ws.Range("E48").Value = '  -0.11%  '

$ws.Range("E49").Value = '  -0.73%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '130.20'
$ws.Range("E50").Value = '  +0.59%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.70'
$ws.Range("E51").Value = '  +1.67%  '
